$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Host all fish images locally: replace remote/mismatched picture paths with local asset paths.
$ws.Range("C11").Value = "assets/fish/rainbow_smelt.jpg"
$ws.Range("C8").Value = "assets/fish/lumpfish.jpeg"

# C11 previously carried the "Hyperlink" cell style (blue/underline) because it held a URL;
# now that it's a plain local path, clear that formatting back to the default style.
$ws.Range("C11").Style = "Normal"

# Reflect the author's final cursor position.
$ws.Range("C9").Select()
